$d = $word.ActiveDocument

# 1. Update the text of the first paragraph: "This is the first level" -> "This is the first level - 1"
$d.Content.Find.Execute("This is the first level", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is the first level - 1", 2)

# 2. Insert a new paragraph right after the (now updated) first paragraph, cloning its
#    numbering/paragraph formatting, with the text "This is the first level - 2".
$firstPara = $d.Paragraphs.First
$firstPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)
$newPara.Range.Text = "This is the first level - 2"

# 3. Remove the trailing duplicate third-level paragraph ("This is the third level - two").
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Delete()
